$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Poland Ekstraklasa")

# --- Block 1: rows 148, 150-154 updated odds/result data (fixtures re-synced to rows) ---
$updates = @(
    @{ Addr = 'B148'; Value = 5465446 }
    @{ Addr = 'F148'; Value = 'Cracovia Krakow' }
    @{ Addr = 'G148'; Value = 'Wisla Plock' }
    @{ Addr = 'H148'; Value = 3 }
    @{ Addr = 'I148'; Value = 0 }
    @{ Addr = 'J148'; Value = 'H' }
    @{ Addr = 'K148'; Value = 2.15 }
    @{ Addr = 'L148'; Value = 3.5 }
    @{ Addr = 'M148'; Value = 2.875 }
    @{ Addr = 'N148'; Value = 2.25 }
    @{ Addr = 'O148'; Value = 3.6 }
    @{ Addr = 'P148'; Value = 2.7 }
    @{ Addr = 'Q148'; Value = -0.25 }
    @{ Addr = 'R148'; Value = 2.05 }
    @{ Addr = 'S148'; Value = 1.75 }
    @{ Addr = 'T148'; Value = 2.5 }
    @{ Addr = 'U148'; Value = 1.825 }
    @{ Addr = 'V148'; Value = 2.025 }
    @{ Addr = 'W148'; Value = 1.25 }
    @{ Addr = 'X148'; Value = -1 }
    @{ Addr = 'Y148'; Value = -1 }
    @{ Addr = 'Z148'; Value = 1.05 }
    @{ Addr = 'AA148'; Value = -1 }
    @{ Addr = 'AB148'; Value = 0.825 }
    @{ Addr = 'AC148'; Value = -1 }
    @{ Addr = 'B150'; Value = 5461474 }
    @{ Addr = 'F150'; Value = 'Legia Warsaw' }
    @{ Addr = 'G150'; Value = 'Slask Wroclaw' }
    @{ Addr = 'H150'; Value = 3 }
    @{ Addr = 'I150'; Value = 1 }
    @{ Addr = 'J150'; Value = 'H' }
    @{ Addr = 'K150'; Value = 1.7 }
    @{ Addr = 'L150'; Value = 3.8 }
    @{ Addr = 'M150'; Value = 4 }
    @{ Addr = 'N150'; Value = 1.833 }
    @{ Addr = 'O150'; Value = 3.8 }
    @{ Addr = 'P150'; Value = 3.4 }
    @{ Addr = 'Q150'; Value = -0.5 }
    @{ Addr = 'R150'; Value = 1.825 }
    @{ Addr = 'S150'; Value = 2.025 }
    @{ Addr = 'T150'; Value = 2.75 }
    @{ Addr = 'U150'; Value = 1.9 }
    @{ Addr = 'V150'; Value = 1.95 }
    @{ Addr = 'W150'; Value = 0.833 }
    @{ Addr = 'X150'; Value = -1 }
    @{ Addr = 'Y150'; Value = -1 }
    @{ Addr = 'Z150'; Value = 0.825 }
    @{ Addr = 'AA150'; Value = -1 }
    @{ Addr = 'AB150'; Value = 0.8999999999999999 }
    @{ Addr = 'AC150'; Value = -1 }
    @{ Addr = 'B151'; Value = 5467427 }
    @{ Addr = 'F151'; Value = 'Stal Mielec' }
    @{ Addr = 'G151'; Value = 'Warta Poznan' }
    @{ Addr = 'H151'; Value = 1 }
    @{ Addr = 'I151'; Value = 0 }
    @{ Addr = 'J151'; Value = 'H' }
    @{ Addr = 'K151'; Value = 2.375 }
    @{ Addr = 'L151'; Value = 3.2 }
    @{ Addr = 'M151'; Value = 2.8 }
    @{ Addr = 'N151'; Value = 2.6 }
    @{ Addr = 'O151'; Value = 3.1 }
    @{ Addr = 'P151'; Value = 2.625 }
    @{ Addr = 'Q151'; Value = 0 }
    @{ Addr = 'R151'; Value = 1.925 }
    @{ Addr = 'S151'; Value = 1.925 }
    @{ Addr = 'T151'; Value = 2.25 }
    @{ Addr = 'U151'; Value = 1.975 }
    @{ Addr = 'V151'; Value = 1.875 }
    @{ Addr = 'W151'; Value = 1.6 }
    @{ Addr = 'X151'; Value = -1 }
    @{ Addr = 'Y151'; Value = -1 }
    @{ Addr = 'Z151'; Value = 0.925 }
    @{ Addr = 'AA151'; Value = -1 }
    @{ Addr = 'AB151'; Value = -1 }
    @{ Addr = 'AC151'; Value = 0.875 }
    @{ Addr = 'B152'; Value = 5456603 }
    @{ Addr = 'F152'; Value = 'Lech Poznan' }
    @{ Addr = 'G152'; Value = 'Jagiellonia Bialystok' }
    @{ Addr = 'H152'; Value = 2 }
    @{ Addr = 'I152'; Value = 0 }
    @{ Addr = 'J152'; Value = 'H' }
    @{ Addr = 'K152'; Value = 1.363 }
    @{ Addr = 'L152'; Value = 4.75 }
    @{ Addr = 'M152'; Value = 6.5 }
    @{ Addr = 'N152'; Value = 1.222 }
    @{ Addr = 'O152'; Value = 5.5 }
    @{ Addr = 'P152'; Value = 8 }
    @{ Addr = 'Q152'; Value = -1.75 }
    @{ Addr = 'R152'; Value = 1.925 }
    @{ Addr = 'S152'; Value = 1.925 }
    @{ Addr = 'T152'; Value = 3.25 }
    @{ Addr = 'U152'; Value = 1.95 }
    @{ Addr = 'V152'; Value = 1.9 }
    @{ Addr = 'W152'; Value = 0.222 }
    @{ Addr = 'X152'; Value = -1 }
    @{ Addr = 'Y152'; Value = -1 }
    @{ Addr = 'Z152'; Value = 0.4625 }
    @{ Addr = 'AA152'; Value = -0.5 }
    @{ Addr = 'AB152'; Value = -1 }
    @{ Addr = 'AC152'; Value = 0.8999999999999999 }
    @{ Addr = 'B153'; Value = 5456594 }
    @{ Addr = 'F153'; Value = 'Rakow Czestochowa' }
    @{ Addr = 'G153'; Value = 'Zaglebie Lubin' }
    @{ Addr = 'H153'; Value = 1 }
    @{ Addr = 'I153'; Value = 1 }
    @{ Addr = 'J153'; Value = 'D' }
    @{ Addr = 'K153'; Value = 1.444 }
    @{ Addr = 'L153'; Value = 4.5 }
    @{ Addr = 'M153'; Value = 5.75 }
    @{ Addr = 'N153'; Value = 1.3 }
    @{ Addr = 'O153'; Value = 5.25 }
    @{ Addr = 'P153'; Value = 7 }
    @{ Addr = 'Q153'; Value = -1.5 }
    @{ Addr = 'R153'; Value = 1.9 }
    @{ Addr = 'S153'; Value = 1.95 }
    @{ Addr = 'T153'; Value = 3 }
    @{ Addr = 'U153'; Value = 1.9 }
    @{ Addr = 'V153'; Value = 1.95 }
    @{ Addr = 'W153'; Value = -1 }
    @{ Addr = 'X153'; Value = 4.25 }
    @{ Addr = 'Y153'; Value = -1 }
    @{ Addr = 'Z153'; Value = -1 }
    @{ Addr = 'AA153'; Value = 0.95 }
    @{ Addr = 'AB153'; Value = -1 }
    @{ Addr = 'AC153'; Value = 0.95 }
    @{ Addr = 'B154'; Value = 5428774 }
    @{ Addr = 'F154'; Value = 'Pogon Szczecin' }
    @{ Addr = 'G154'; Value = 'Radomiak Radom' }
    @{ Addr = 'H154'; Value = 4 }
    @{ Addr = 'I154'; Value = 0 }
    @{ Addr = 'J154'; Value = 'H' }
    @{ Addr = 'K154'; Value = 1.571 }
    @{ Addr = 'L154'; Value = 4 }
    @{ Addr = 'M154'; Value = 4.75 }
    @{ Addr = 'N154'; Value = 1.533 }
    @{ Addr = 'O154'; Value = 4.333 }
    @{ Addr = 'P154'; Value = 4.75 }
    @{ Addr = 'Q154'; Value = -1 }
    @{ Addr = 'R154'; Value = 1.875 }
    @{ Addr = 'S154'; Value = 1.975 }
    @{ Addr = 'T154'; Value = 3 }
    @{ Addr = 'U154'; Value = 1.875 }
    @{ Addr = 'V154'; Value = 1.975 }
    @{ Addr = 'W154'; Value = 0.5329999999999999 }
    @{ Addr = 'X154'; Value = -1 }
    @{ Addr = 'Y154'; Value = -1 }
    @{ Addr = 'Z154'; Value = 0.875 }
    @{ Addr = 'AA154'; Value = -1 }
    @{ Addr = 'AB154'; Value = 0.875 }
    @{ Addr = 'AC154'; Value = -1 }

# --- Block 2: rows 341-349 refreshed closing-odds values for upcoming fixtures ---
    @{ Addr = 'O341'; Value = 2.9 }
    @{ Addr = 'Q341'; Value = -0.25 }
    @{ Addr = 'R341'; Value = 1.775 }
    @{ Addr = 'S341'; Value = 2.1 }
    @{ Addr = 'U341'; Value = 2.025 }
    @{ Addr = 'V341'; Value = 1.825 }
    @{ Addr = 'N342'; Value = 1.3 }
    @{ Addr = 'O342'; Value = 5.25 }
    @{ Addr = 'P342'; Value = 8.5 }
    @{ Addr = 'R342'; Value = 1.875 }
    @{ Addr = 'S342'; Value = 1.975 }
    @{ Addr = 'U342'; Value = 1.825 }
    @{ Addr = 'V342'; Value = 2.025 }
    @{ Addr = 'U343'; Value = 1.875 }
    @{ Addr = 'V343'; Value = 1.975 }
    @{ Addr = 'R344'; Value = 1.8 }
    @{ Addr = 'S344'; Value = 2.05 }
    @{ Addr = 'U344'; Value = 1.9 }
    @{ Addr = 'V344'; Value = 1.95 }
    @{ Addr = 'N345'; Value = 5 }
    @{ Addr = 'O345'; Value = 3.8 }
    @{ Addr = 'P345'; Value = 1.65 }
    @{ Addr = 'R345'; Value = 2.025 }
    @{ Addr = 'S345'; Value = 1.825 }
    @{ Addr = 'U345'; Value = 1.825 }
    @{ Addr = 'V345'; Value = 2.025 }
    @{ Addr = 'R346'; Value = 1.925 }
    @{ Addr = 'S346'; Value = 1.925 }
    @{ Addr = 'R349'; Value = 1.8 }
    @{ Addr = 'S349'; Value = 2.05 }
)

foreach ($item in $updates) {
    $ws.Range($item.Addr).Value = $item.Value
}
